$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Set Runmode (column D) to "Y" for all test case rows (2 through 18)
$ws.Range("D2:D18").Value = "Y"

# The test in row 6 has not actually been run yet, so its Result goes back to SKIP
$ws.Range("E6").Value = "SKIP"

# Reflect the selection left behind after selecting D2:D18 to run all the tests
$ws.Range("D2:D18").Select()
